$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 59, shifting the existing
# rows 59-71 down to 61-73 (a new week of prices is prepended to the list).
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(60).Insert()

# New row 59: Primera quality, week of 2021-10-22 (serial 44491)
$ws.Range("A59").Value = 1
$ws.Range("B59").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C59").Value = "Arica y Parinacota"
$ws.Range("D59").Value = 44491
$ws.Range("E59").Value = 15
$ws.Range("F59").Value = 100112036
$ws.Range("G59").Value = "Caigua"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 100
$ws.Range("K59").Value = 7000
$ws.Range("L59").Value = 8000
$ws.Range("M59").Value = 7500
$ws.Range("N59").Value = "`$/caja 20 kilos"
$ws.Range("O59").Value = "Región de Arica y Parinacota"
$ws.Range("P59").Value = 375
$ws.Range("Q59").Value = 20
$ws.Range("R59").Value = "Hortaliza"

# New row 60: Segunda quality, same week
$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C60").Value = "Arica y Parinacota"
$ws.Range("D60").Value = 44491
$ws.Range("E60").Value = 15
$ws.Range("F60").Value = 100112036
$ws.Range("G60").Value = "Caigua"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Segunda"
$ws.Range("J60").Value = 120
$ws.Range("K60").Value = 5000
$ws.Range("L60").Value = 6000
$ws.Range("M60").Value = 5500
$ws.Range("N60").Value = "`$/caja 20 kilos"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 275
$ws.Range("Q60").Value = 20
$ws.Range("R60").Value = "Hortaliza"

# Append two new rows 72-73 at the end (copies of the former last two
# rows, which have now shifted to 70 and 71 after the insert above).
$ws.Range("A72").Value = 1
$ws.Range("B72").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C72").Value = "Arica y Parinacota"
$ws.Range("D72").Value = 44438
$ws.Range("E72").Value = 15
$ws.Range("F72").Value = 100112036
$ws.Range("G72").Value = "Caigua"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 120
$ws.Range("K72").Value = 6000
$ws.Range("L72").Value = 7000
$ws.Range("M72").Value = 6500
$ws.Range("N72").Value = "`$/caja 20 kilos"
$ws.Range("O72").Value = "Región de Arica y Parinacota"
$ws.Range("P72").Value = 325
$ws.Range("Q72").Value = 20
$ws.Range("R72").Value = "Hortaliza"

$ws.Range("A73").Value = 1
$ws.Range("B73").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C73").Value = "Arica y Parinacota"
$ws.Range("D73").Value = 44438
$ws.Range("E73").Value = 15
$ws.Range("F73").Value = 100112036
$ws.Range("G73").Value = "Caigua"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Segunda"
$ws.Range("J73").Value = 120
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = 5500
$ws.Range("N73").Value = "`$/caja 20 kilos"
$ws.Range("O73").Value = "Región de Arica y Parinacota"
$ws.Range("P73").Value = 275
$ws.Range("Q73").Value = 20
$ws.Range("R73").Value = "Hortaliza"

# Apply the date number format used by the rest of column D to the two
# newly appended rows, matching the existing column D formatting.
$ws.Range("D72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D73").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Final UsedRange:" $ws.UsedRange.Address()
